$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C/D to D/E), mirroring the
# "right-click column C header -> Insert" workflow used to add the VIP column.
$ws.Columns("C").Insert()

# Seed the new VIP column with the BC column's values as a starting point,
# then rename the header from "BC" to "VIP".
$ws.Range("B1:B31").Copy($ws.Range("C1:C31"))
$ws.Range("C1").Value = "VIP"

# Match the column's on-disk width (best achievable precision through the
# pixel-quantized ColumnWidth COM setter).
$ws.Columns("C").ColumnWidth = 16.5

# Leave the whole new column selected, as if freshly inserted via the UI.
$ws.Columns("C").Select()
